$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows down
$ws.Rows("2:2").Insert()
$ws.Rows("2:2").ClearFormats()

# Populate the new row 2 with the new weekly data
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C2").Value = "Los Lagos"
$ws.Range("D2").Value = 44530
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100101
$ws.Range("H2").Value = "Berries"
$ws.Range("I2").Value = 100101001
$ws.Range("J2").Value = "Arándano (blue)"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 3600
$ws.Range("O2").Value = 3700
$ws.Range("P2").Value = 3650
$ws.Range("Q2").Value = "$/kilo"
$ws.Range("R2").Value = "Región del Maule"
$ws.Range("S2").Value = 3650
$ws.Range("T2").Value = 1

# Match the date number format used by the other rows in column D
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat
